$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.97%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.97%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.470"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08076"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.76%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.663"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.69%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.94%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.882"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.57%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9376"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.40%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1197"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.06%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1895"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.39%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09678"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.56%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04077"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "4.78%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1068"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.63%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001273"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.27%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005979"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.17%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.578"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.46%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.04%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.610"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.53%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1332"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.26%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2497"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.55%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04347"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.16%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001236"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.76%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004276"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.49%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001234"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.59%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004007"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.33%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02655"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.75%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05451"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.33%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007664"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.63%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009726"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.28%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.44%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002128"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.66%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009905"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-15.56%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007131"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.65%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.33%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003573"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "3.25%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002280"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.00%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002110"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.33%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002009"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.33%"
